$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Start from a clean sheet so the row/column layout below is authoritative.
# ---------------------------------------------------------------------------
$ws.Cells.Clear()

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.7109375
$ws.Columns.Item(2).ColumnWidth = 20.7109375
$ws.Columns.Item(3).ColumnWidth = 15.7109375
$ws.Columns.Item(4).ColumnWidth = 15.7109375

# ---------------------------------------------------------------------------
# Header row (bold + centered)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Ranges"
$ws.Range("B1").Value = "Old"
$ws.Range("C1").Value = "Current"
$ws.Range("D1").Value = "Old % New"
$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Per Year block
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Per Year Min:"
$ws.Range("B2").Value = 200000
$ws.Range("C2").Value = 465000
$ws.Range("D2").Formula = "=(B2/C2)*100"

$ws.Range("A3").Value = "Per Year Max:"
$ws.Range("B3").Value = 14171500
$ws.Range("C3").Value = 23943600
$ws.Range("D3").Formula = "=(B3/C3)*100"

# ---------------------------------------------------------------------------
# NoBon Per Year block
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "NoBon Per Year Min:"
$ws.Range("B5").Value = 10000
$ws.Range("C5").Value = 41000
$ws.Range("D5").Formula = "=(B5/C5)*100"

$ws.Range("A6").Value = "NoBon Per Year Max:"
$ws.Range("B6").Value = 10811500
$ws.Range("C6").Value = 14170833
$ws.Range("D6").Formula = "=(B6/C6)*100"

# ---------------------------------------------------------------------------
# Total Salary block
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Total Salary Min:"
$ws.Range("B8").Value = 230000
$ws.Range("C8").Value = 465000
$ws.Range("D8").Formula = "=(B8/C8)*100"

$ws.Range("A9").Value = "Total Salary Max:"
$ws.Range("B9").Value = 99200000
$ws.Range("C9").Value = 135000000
$ws.Range("D9").Formula = "=(B9/C9)*100"

# ---------------------------------------------------------------------------
# Signing Bonus block
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "Signing Bonus Min:"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "Signing Bonus Max:"
$ws.Range("B12").Value = 23520000
$ws.Range("C12").Value = 60500000
$ws.Range("D12").Formula = "=(B12/C12)*100"

# Blank spacer rows (4, 7, 10) still get the numeric style applied to B/C,
# matching the look of the populated rows around them.
$ws.Range("B4:C4,B7:C7,B10:C10").Value = ""

# Apply the "#,##0" number style to every B/C cell in the ranges table
# (rows 2-12), including the blank spacer rows.
$ws.Range("B2:C12").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# Proposed Rules section
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Proposed Rules:"

$ws.Range("A16").Value = "Total Sal < 1M, cut to 43%"
$ws.Range("A17").Value = "Total Sal > 1M, < 10M, cut to 58%"
$ws.Range("A18").Value = "Total Sal > 10M, cut to 72.5%"

$ws.Range("A19").Value = "Formula:"
$ws.Range("B19").Value = "if tot_sal > 10M: "
$ws.Range("C19").Value = "new_tot_sal = round((tot_sal / 10000) * .725)"

$ws.Range("B20").Value = "elif tot_sal > 1M:"
$ws.Range("C20").Value = "new_tot_sal = round((tot_sal / 10000) * .58)"

$ws.Range("B21").Value = "else:"
$ws.Range("C21").Value = "new_tot_sal = round((tot_sal / 10000) * .43)"

$ws.Range("A23").Value = "NOTE!! PSBO must always be in multiples of PCON (length of contract). This means if the contract is for 3 years, the bonus can NOT be 10,000 or 20,000 (PSBO != 1 or 2)."

$ws.Range("A24").Value = "Signing Bon < 100K, cut to 80%"
$ws.Range("A25").Value = "Signing Bon >100K, < 1M, cut to 65%"
$ws.Range("A26").Value = "Signing Bon >1M, < 10M, cut to 50%"
$ws.Range("A27").Value = "Signing Bon >10M, cut to 40%"

$ws.Range("A28").Value = "Formula:"
$ws.Range("B28").Value = "if sgn_bon > 10M: "
$ws.Range("C28").Value = "new_sgn_bon = round((sgn_bon / 10000) * .4)"

$ws.Range("B29").Value = "elif sgn_bon > 1M: "
$ws.Range("C29").Value = "new_sgn_bon = round((sgn_bon / 10000) * .5)"

$ws.Range("B30").Value = "elif sgn_bon > 100K: "
$ws.Range("C30").Value = "new_sgn_bon = round((sgn_bon / 10000) * .65)"

$ws.Range("B31").Value = "else:"
$ws.Range("C31").Value = "new_sgn_bon = round((sgn_bon / 10000) * .8)"

$ws.Range("B32").Value = "AND then:"
$ws.Range("C32").Value = "if (new_sgn_bon Mod PCON) > 0: new_sgn_bon +=  PCON - (new_sgn_bon Mod PCON)"

# ---------------------------------------------------------------------------
# Page setup
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

Write-Host "Edit complete"
